# Update the YouTube link URL in the "Coding Assignment Week 13" document.
# The original run reads:
#   "YouTube Link: https://youtu.be/8_d240eXPME"
# and must become two runs (same bold formatting):
#   "YouTube Link: " + "https://youtu.be/56vsw1zcX6A"

$d = $word.ActiveDocument

$oldUrl = "https://youtu.be/8_d240eXPME"
$newUrl = "https://youtu.be/56vsw1zcX6A"

# Replace just the URL portion, leaving the "YouTube Link: " label untouched.
$find = $d.Content
$found = $find.Find.Execute($oldUrl, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $newUrl, 2)

if (-not $found) {
    throw "Could not find the YouTube URL to replace."
}

# Locate the freshly-inserted URL text and toggle Bold off/on so that it
# becomes a distinct run from the preceding "YouTube Link: " label (Word
# otherwise keeps them merged into a single run since the formatting is
# identical).
$urlRange = $d.Content
$urlFound = $urlRange.Find.Execute($newUrl)

if (-not $urlFound) {
    throw "Could not locate the newly inserted YouTube URL."
}

$urlRange.Font.Bold = 0
$urlRange.Font.Bold = 1
